$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '39.803.60'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.200.88'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.67%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '291.23'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.56%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '86.22'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.508'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.66%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.466'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.02%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '30.10'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.56%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '50.02'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +6.37%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0776'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.78%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.112'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.63%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.43'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.538.53'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.76%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '13.68'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -3.66%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.154.60'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -4.20%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.727'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '39.724.83'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.93%  '
$ws.Range('E20').Value = '  -0.89%  '
$ws.Range('E21').Value = '  -1.47%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.72'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.00%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '65.15'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.07%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '237.32'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.44'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.80%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.81'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.90%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '23.44'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.11%  '
$ws.Range('E29').Value = '  -7.63%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.17'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.65%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '155.96'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.70%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '31.24'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -6.12%  '
$ws.Range('E33').Value = '  -0.10%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.92'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.25%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0705'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.52%  '
$ws.Range('E36').Value = '  -2.45%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.86'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.23%  '
$ws.Range('E38').Value = '  -0.31%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0974'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.75%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '15.18'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -6.98%  '
$ws.Range('E41').Value = '  -2.43%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.111.84'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.11%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.72'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.68%  '
$ws.Range('E44').Value = '  -0.74%  '
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '9.77'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.34%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '17.38'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.13%  '
$ws.Range('E48').Value = '  +1.63%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.409.61'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.87%  '
$ws.Range('E50').Value = '  +0.64%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '88.23'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.43%  '
